$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$jValues = @{
    2 = 4
    3 = 5
    4 = 5
    5 = 3
    6 = 4
    7 = 5
    8 = 5
    9 = 2
    10 = 1
    11 = 1
    12 = 1
    13 = 2
    14 = 2
    15 = 5
    16 = 2
    17 = 4
    18 = 3
    19 = 1
    20 = 2
    21 = 4
    22 = 3
    23 = 4
    24 = 4
    25 = 5
    26 = 3
    27 = 1
    28 = 2
    29 = 2
    30 = 5
    31 = 4
    32 = 3
    33 = 4
    34 = 5
    35 = 1
    36 = 2
    37 = 1
    38 = 3
    39 = 5
    40 = 4
    41 = 5
    42 = 2
    43 = 5
    44 = 4
    45 = 2
    46 = 3
    47 = 2
    48 = 2
    49 = 3
    50 = 4
    52 = 4
    53 = 2
    54 = 5
    55 = 5
    56 = 4
    57 = 3
    58 = 4
    59 = 1
    60 = 3
    61 = 5
    62 = 5
    63 = 3
    64 = 4
    65 = 5
    66 = 5
    67 = 2
    68 = 2
    69 = 2
    70 = 1
    71 = 1
    72 = 1
    73 = 3
    74 = 3
    75 = 2
    76 = 2
    77 = 4
    78 = 4
    79 = 1
    80 = 3
    81 = 1
    82 = 5
    83 = 2
    84 = 5
    85 = 4
    86 = 2
    87 = 4
    88 = 1
    89 = 5
    90 = 5
    91 = 2
    92 = 2
    93 = 5
    94 = 3
    95 = 1
    96 = 5
    97 = 3
    98 = 4
    99 = 1
    100 = 5
}

foreach ($row in $jValues.Keys) {
    $ws.Cells.Item([int]$row, 10).Value = $jValues[$row]
}
